$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.988.77"
$ws.Range("E2").Value = "'  +4.42%  "
$ws.Range("D3").Value = "'3.244.30"
$ws.Range("E3").Value = "'  +2.18%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'396.27"
$ws.Range("E5").Value = "'  -1.23%  "
$ws.Range("D6").Value = "'108.31"
$ws.Range("E6").Value = "'  -1.47%  "
$ws.Range("E7").Value = "'  +7.14%  "
$ws.Range("D8").Value = "'3.238.85"
$ws.Range("E8").Value = "'  +2.20%  "
$ws.Range("E9").Value = "'  +0.02%  "
$ws.Range("E10").Value = "'  +1.46%  "
$ws.Range("D11").Value = "'39.27"
$ws.Range("E11").Value = "'  +0.40%  "
$ws.Range("D12").Value = "'0.0991"
$ws.Range("E12").Value = "'  +10.62%  "
$ws.Range("E13").Value = "'  +2.10%  "
$ws.Range("D14").Value = "'3.755.67"
$ws.Range("D15").Value = "'8.34"
$ws.Range("E15").Value = "'  +3.30%  "
$ws.Range("D16").Value = "'19.10"
$ws.Range("E16").Value = "'  +0.21%  "
$ws.Range("D17").Value = "'3.243.92"
$ws.Range("E17").Value = "'  +2.30%  "
$ws.Range("E18").Value = "'  -2.82%  "
$ws.Range("D19").Value = "'10.75"
$ws.Range("E19").Value = "'  +1.92%  "
$ws.Range("D20").Value = "'56.823.50"
$ws.Range("E20").Value = "'  +4.30%  "
$ws.Range("D21").Value = "'3.34"
$ws.Range("E21").Value = "'  +1.32%  "
$ws.Range("E22").Value = "'  +7.41%  "
$ws.Range("E23").Value = "'  +0.92%  "
$ws.Range("D24").Value = "'294.95"
$ws.Range("E24").Value = "'  +6.91%  "
$ws.Range("D25").Value = "'74.23"
$ws.Range("E25").Value = "'  +2.79%  "
$ws.Range("E26").Value = "'  -2.37%  "
$ws.Range("D27").Value = "'28.09"
$ws.Range("E27").Value = "'  +1.10%  "
$ws.Range("E28").Value = "'  +1.04%  "
$ws.Range("D29").Value = "'7.65"
$ws.Range("E29").Value = "'  -5.41%  "
$ws.Range("D30").Value = "'7.24"
$ws.Range("E30").Value = "'  -4.92%  "
$ws.Range("E31").Value = "'  -0.76%  "
$ws.Range("E32").Value = "'  -0.02%  "
$ws.Range("D33").Value = "'11.22"
$ws.Range("E33").Value = "'  +1.46%  "
$ws.Range("E34").Value = "'  -3.05%  "
$ws.Range("D35").Value = "'39.99"
$ws.Range("D36").Value = "'0.0489"
$ws.Range("E36").Value = "'  -3.42%  "
$ws.Range("E37").Value = "'  +1.26%  "
$ws.Range("D38").Value = "'51.50"
$ws.Range("E38").Value = "'  +0.13%  "
$ws.Range("D39").Value = "'0.998"
$ws.Range("E39").Value = "'  -0.15%  "
$ws.Range("D40").Value = "'3.48"
$ws.Range("E40").Value = "'  -4.30%  "
$ws.Range("E41").Value = "'  +1.38%  "
$ws.Range("D42").Value = "'138.92"
$ws.Range("E42").Value = "'  +5.52%  "
$ws.Range("E43").Value = "'  +4.01%  "
$ws.Range("E44").Value = "'  -1.61%  "
$ws.Range("D45").Value = "'17.07"
$ws.Range("E45").Value = "'  -0.80%  "
$ws.Range("E46").Value = "'  -3.44%  "
$ws.Range("E47").Value = "'  -4.00%  "
$ws.Range("D48").Value = "'22.17"
$ws.Range("E48").Value = "'  +0.26%  "
$ws.Range("D49").Value = "'2.19"
$ws.Range("E49").Value = "'  +5.56%  "
$ws.Range("D50").Value = "'2.161.80"
$ws.Range("E50").Value = "'  +2.94%  "
$ws.Range("E51").Value = "'  -6.33%  "
